$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 3299.6667
$ws.Range("J38").Value = 2849.5
$ws.Range("L38").Value = 8548.5
$ws.Range("N38").Value = -9292.5

# Row 58
$ws.Range("H58").Value = 5309.591
$ws.Range("I58").Value = 2100.1538
$ws.Range("J58").Value = 9945.444
$ws.Range("K58").Value = 6300.4614
$ws.Range("L58").Value = 29836.332
$ws.Range("M58").Value = -6150.4614
$ws.Range("N58").Value = -30136.332

# Row 62
$ws.Range("H62").Value = 6948533
$ws.Range("I62").Value = 11366329
$ws.Range("K62").Value = 11366329
$ws.Range("M62").Value = -11365705

# Row 65
$ws.Range("H65").Value = 6948533
$ws.Range("I65").Value = 11366329
$ws.Range("K65").Value = 56831645
$ws.Range("M65").Value = -56828525

# Row 96
$ws.Range("H96").Value = 704.9091
$ws.Range("J96").Value = 999.5
$ws.Range("L96").Value = 2998.5
$ws.Range("N96").Value = -5744.5

# Row 99
$ws.Range("H99").Value = 590.9231
$ws.Range("I99").Value = 264.66666
$ws.Range("J99").Value = 870.5714
$ws.Range("K99").Value = 793.9999799999999
$ws.Range("L99").Value = 2611.7142
$ws.Range("M99").Value = 704.0000200000001
$ws.Range("N99").Value = -5607.7142

# Row 106
$ws.Range("H106").Value = 2799.7856
$ws.Range("I106").Value = 2630.923
$ws.Range("K106").Value = 2630.923
$ws.Range("M106").Value = -1999.923

# Row 131
$ws.Range("H131").Value = 3636.739
$ws.Range("I131").Value = 1727.8125
$ws.Range("J131").Value = 8000
$ws.Range("K131").Value = 5183.4375
$ws.Range("L131").Value = 24000
$ws.Range("M131").Value = -143.4375
$ws.Range("N131").Value = -34080

# Row 132
$ws.Range("H132").Value = 1976.175
$ws.Range("I132").Value = 1834.5758
$ws.Range("J132").Value = 2643.7144
$ws.Range("K132").Value = 5503.7274
$ws.Range("L132").Value = 7931.1432
$ws.Range("M132").Value = -2973.7274
$ws.Range("N132").Value = -12991.1432

# Row 137
$ws.Range("H137").Value = 2488.9473
$ws.Range("I137").Value = 2046.4706
$ws.Range("K137").Value = 6139.4118
$ws.Range("M137").Value = -3589.4118

$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Range("H24").Value = 17499.5

# Row 32
$ws.Range("H32").Value = 2289.34
$ws.Range("I32").Value = 2261.9597
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 2261.9597
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = -1974.9597
$ws.Range("N32").Value = -5574

# Row 45
$ws.Range("H45").Value = 2154.6667
$ws.Range("I45").Value = 2478.4
$ws.Range("K45").Value = 2478.4
$ws.Range("M45").Value = -2101.4

# Row 74
$ws.Range("H74").Value = 1950.1
$ws.Range("I74").Value = 1966.8
$ws.Range("J74").Value = 1900
$ws.Range("K74").Value = 1966.8
$ws.Range("L74").Value = 1900
$ws.Range("M74").Value = -1092.8
$ws.Range("N74").Value = -3648

# Row 77
$ws.Range("H77").Value = 1950.1
$ws.Range("I77").Value = 1966.8
$ws.Range("J77").Value = 1900
$ws.Range("K77").Value = 9834
$ws.Range("L77").Value = 9500
$ws.Range("M77").Value = -5466
$ws.Range("N77").Value = -18236

# Row 97
$ws.Range("H97").Value = 959.2083
$ws.Range("I97").Value = 993.3043
$ws.Range("J97").Value = 175
$ws.Range("K97").Value = 993.3043
$ws.Range("L97").Value = 175
$ws.Range("M97").Value = -497.3043
$ws.Range("N97").Value = -1167

# Row 100
$ws.Range("H100").Value = 17499.5

# Row 106
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("N106").ClearContents()

# Row 122
$ws.Range("H122").Value = 4668.9585
$ws.Range("I122").Value = 3191.3635
$ws.Range("K122").Value = 9574.0905
$ws.Range("M122").Value = -7124.0905

# Row 132
$ws.Range("H132").Value = 2816.3044
$ws.Range("I132").Value = 2751.262
$ws.Range("K132").Value = 8253.786
$ws.Range("M132").Value = -5723.786

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3091.2354
$ws.Range("I20").Value = 2776.2083
$ws.Range("J20").Value = 3847.3
$ws.Range("K20").Value = 2776.2083
$ws.Range("L20").Value = 3847.3
$ws.Range("M20").Value = -2529.2083
$ws.Range("N20").Value = -4341.3

# Row 134
$ws.Range("H134").Value = 57193.25
$ws.Range("I134").Value = 7436.9443
$ws.Range("K134").Value = 22310.8329
$ws.Range("M134").Value = -19775.8329

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 1036.5714
$ws.Range("I23").Value = 287.66666
$ws.Range("K23").Value = 862.9999799999999
$ws.Range("M23").Value = -627.9999799999999

# Row 80
$ws.Range("H80").Value = 2000.6
$ws.Range("J80").Value = 2333.6667
$ws.Range("L80").Value = 7001.000100000001
$ws.Range("N80").Value = -8873.000100000001

# Row 83
$ws.Range("H83").Value = 2000.6
$ws.Range("J83").Value = 2333.6667
$ws.Range("L83").Value = 21003.0003
$ws.Range("N83").Value = -30363.0003

# Row 87
$ws.Range("H87").Value = 18998
$ws.Range("I87").Value = 18998
$ws.Range("K87").Value = 56994
$ws.Range("M87").Value = -55746

# Row 90
$ws.Range("H90").Value = 18998
$ws.Range("I90").Value = 18998
$ws.Range("K90").Value = 170982
$ws.Range("M90").Value = -164742

# Row 92
$ws.Range("H92").Value = 910109.2
$ws.Range("I92").Value = 1250434.9
$ws.Range("J92").Value = 2574
$ws.Range("K92").Value = 3751304.7
$ws.Range("L92").Value = 7722
$ws.Range("M92").Value = -3750056.7
$ws.Range("N92").Value = -10218

# Row 104
$ws.Range("H104").Value = 4626
$ws.Range("I104").Value = 4626
$ws.Range("K104").Value = 13878
$ws.Range("M104").Value = -11257

# Row 107
$ws.Range("H107").Value = 76062.14
$ws.Range("I107").Value = 1146.2858
$ws.Range("J107").Value = 150978
$ws.Range("K107").Value = 3438.8574
$ws.Range("L107").Value = 452934
$ws.Range("M107").Value = -1518.8574
$ws.Range("N107").Value = -456774

# Row 109
$ws.Range("H109").Value = 46874.824
$ws.Range("I109").Value = 1481.2222
$ws.Range("K109").Value = 4443.6666
$ws.Range("M109").Value = -3403.6666

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 2374.25
$ws.Range("J122").Value = 2749.5
$ws.Range("L122").Value = 8248.5
$ws.Range("N122").Value = -13148.5

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 2628.1428
$ws.Range("I68").Value = 2564.75
$ws.Range("K68").Value = 2564.75
$ws.Range("M68").Value = -1815.75

# Row 71
$ws.Range("H71").Value = 2628.1428
$ws.Range("I71").Value = 2564.75
$ws.Range("K71").Value = 12823.75
$ws.Range("M71").Value = -9079.75

# Row 100
$ws.Range("H100").Value = 4936.125
$ws.Range("J100").Value = 4899
$ws.Range("L100").Value = 4899
$ws.Range("N100").Value = -5981

# Row 132
$ws.Range("H132").Value = 4527.619
$ws.Range("J132").Value = 901.5
$ws.Range("L132").Value = 2704.5
$ws.Range("N132").Value = -7764.5

# Row 136
$ws.Range("H136").Value = 351403.72
$ws.Range("I136").Value = 592841.2
$ws.Range("J136").Value = 9367.333000000001
$ws.Range("K136").Value = 1778523.6
$ws.Range("L136").Value = 28101.999
$ws.Range("N136").Value = -33201.999

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 396.66666
$ws.Range("I113").Value = 321.25
$ws.Range("K113").Value = 963.75
$ws.Range("M113").Value = 1206.25

# Row 132
$ws.Range("H132").Value = 30759.352
$ws.Range("I132").Value = 2873.5557
$ws.Range("J132").Value = 106051
$ws.Range("K132").Value = 8620.667099999999
$ws.Range("L132").Value = 318153
$ws.Range("M132").Value = -6090.667099999999
$ws.Range("N132").Value = -323213

# Row 136
$ws.Range("H136").Value = 120774
$ws.Range("I136").Value = 3322.375
$ws.Range("K136").Value = 9967.125
$ws.Range("M136").Value = -7417.125
